# changing azure route table map for new naming scheme
#
# - Rename the "WANCPFW" worksheet to "CPMGMT" (per the new naming scheme).
# - Make that sheet the active/selected tab (it was previously on Fortinet_1).
# - Move the sheet's selection to E14.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("WANCPFW")
$ws.Name = "CPMGMT"

$ws.Activate()
$ws.Range("E14").Select()
